$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted right after the existing row 105
# (2021-10-07), pushing all subsequent rows (106-127) down by one to
# (107-128). Insert a fresh row at position 106 to achieve that shift.
$ws.Range("A106:R106").EntireRow.Insert()

# Populate the newly inserted row 106 with the new weekly record.
$ws.Range("A106").Value = 1
$ws.Range("B106").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C106").Value = "Arica y Parinacota"
$ws.Range("D106").Value = 44798
$ws.Range("E106").Value = 15
$ws.Range("F106").Value = 100112036
$ws.Range("G106").Value = "Caigua"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 130
$ws.Range("K106").Value = 7000
$ws.Range("L106").Value = 8000
$ws.Range("M106").Value = 7500
$ws.Range("N106").Value = "$/caja 20 kilos"
$ws.Range("O106").Value = "Región de Arica y Parinacota"
$ws.Range("P106").Value = 375
$ws.Range("Q106").Value = 20
$ws.Range("R106").Value = "Hortaliza"
